# Update the hyperlink/photo path text on Sheet3 (shared string used by B2,
# "Photo Add" row) to the new file path.
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B2").Value = "F:\GraphicsViewer.aspx.png"

# Switch the active/selected sheet from Sheet2 to Sheet3.
$ws3.Activate()
